$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.845.81'
$ws.Cells.Item(2, 5).Value = '  -0.22%  '

$ws.Cells.Item(3, 4).Value = '1.735.89'
$ws.Cells.Item(3, 5).Value = '  -0.37%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '240.64'
$ws.Cells.Item(5, 5).Value = '  +4.12%  '

$ws.Cells.Item(6, 5).Value = '  +0.05%  '

$ws.Cells.Item(7, 4).Value = '0.5187'
$ws.Cells.Item(7, 5).Value = '  -1.03%  '

$ws.Cells.Item(8, 4).Value = '0.2736'
$ws.Cells.Item(8, 5).Value = '  -0.93%  '

$ws.Cells.Item(9, 5).Value = '  -0.01%  '

$ws.Cells.Item(10, 4).Value = '1.739.55'
$ws.Cells.Item(10, 5).Value = '  +0.21%  '

$ws.Cells.Item(11, 4).Value = '0.07167'
$ws.Cells.Item(11, 5).Value = '  +0.66%  '

$ws.Cells.Item(12, 4).Value = '14.94'
$ws.Cells.Item(12, 5).Value = '  -2.00%  '

$ws.Cells.Item(13, 4).Value = '0.6404'
$ws.Cells.Item(13, 5).Value = '  -0.93%  '

$ws.Cells.Item(14, 4).Value = '4.602'
$ws.Cells.Item(14, 5).Value = '  +1.70%  '

$ws.Cells.Item(15, 4).Value = '77.09'
$ws.Cells.Item(15, 5).Value = '  -0.31%  '

$ws.Cells.Item(16, 5).Value = '  +0.04%  '

$ws.Cells.Item(17, 5).Value = '  +0.00%  '

$ws.Cells.Item(18, 4).Value = '25.880.16'
$ws.Cells.Item(18, 5).Value = '  -0.06%  '

$ws.Cells.Item(19, 4).Value = '11.71'
$ws.Cells.Item(19, 5).Value = '  +1.42%  '

$ws.Cells.Item(20, 4).Value = '0.000006751'
$ws.Cells.Item(20, 5).Value = '  +0.96%  '

$ws.Cells.Item(21, 4).Value = '1.963.87'
$ws.Cells.Item(21, 5).Value = '  +0.22%  '

$ws.Cells.Item(22, 4).Value = '4.257'
$ws.Cells.Item(22, 5).Value = '  -0.50%  '

$ws.Cells.Item(23, 4).Value = '8.615'
$ws.Cells.Item(23, 5).Value = '  -1.93%  '

$ws.Cells.Item(24, 4).Value = '5.252'
$ws.Cells.Item(24, 5).Value = '  +1.38%  '

$ws.Cells.Item(25, 4).Value = '137.70'
$ws.Cells.Item(25, 5).Value = '  -1.59%  '

$ws.Cells.Item(26, 4).Value = '1.507'
$ws.Cells.Item(26, 5).Value = '  -1.01%  '

$ws.Cells.Item(27, 4).Value = '15.17'
$ws.Cells.Item(27, 5).Value = '  -0.21%  '

$ws.Cells.Item(28, 4).Value = '1.764'
$ws.Cells.Item(28, 5).Value = '  -2.23%  '

$ws.Cells.Item(29, 4).Value = '104.88'
$ws.Cells.Item(29, 5).Value = '  +2.22%  '

$ws.Cells.Item(30, 4).Value = '3.936'
$ws.Cells.Item(30, 5).Value = '  +5.38%  '

$ws.Cells.Item(31, 4).Value = '0.08238'
$ws.Cells.Item(31, 5).Value = '  -1.18%  '

$ws.Cells.Item(32, 4).Value = '3.647'
$ws.Cells.Item(32, 5).Value = '  +1.80%  '

$ws.Cells.Item(33, 4).Value = '0.04626'
$ws.Cells.Item(33, 5).Value = '  +2.21%  '

$ws.Cells.Item(34, 4).Value = '2.644'
$ws.Cells.Item(34, 5).Value = '  +1.23%  '

$ws.Cells.Item(35, 4).Value = '0.9857'
$ws.Cells.Item(35, 5).Value = '  +0.54%  '

$ws.Cells.Item(36, 4).Value = '0.6164'
$ws.Cells.Item(36, 5).Value = '  -0.94%  '

$ws.Cells.Item(37, 4).Value = '2.685'
$ws.Cells.Item(37, 5).Value = '  -0.33%  '

$ws.Cells.Item(38, 5).Value = '  +0.45%  '

$ws.Cells.Item(39, 4).Value = '1.917'
$ws.Cells.Item(39, 5).Value = '  -0.56%  '

$ws.Cells.Item(41, 4).Value = '99.95'
$ws.Cells.Item(41, 5).Value = '  -0.36%  '

$ws.Cells.Item(42, 4).Value = '0.3832'
$ws.Cells.Item(42, 5).Value = '  -1.25%  '

$ws.Cells.Item(43, 4).Value = '0.7447'
$ws.Cells.Item(43, 5).Value = '  +1.43%  '

$ws.Cells.Item(44, 4).Value = '4.983'
$ws.Cells.Item(44, 5).Value = '  -0.79%  '

$ws.Cells.Item(45, 4).Value = '0.1121'
$ws.Cells.Item(45, 5).Value = '  -0.46%  '

$ws.Cells.Item(46, 4).Value = '6.232'
$ws.Cells.Item(46, 5).Value = '  -0.43%  '

$ws.Cells.Item(47, 4).Value = '0.05241'
$ws.Cells.Item(47, 5).Value = '  -1.84%  '

$ws.Cells.Item(48, 4).Value = '54.79'
$ws.Cells.Item(48, 5).Value = '  +1.96%  '

$ws.Cells.Item(49, 4).Value = '30.47'
$ws.Cells.Item(49, 5).Value = '  +0.77%  '

$ws.Cells.Item(50, 4).Value = '7.544'
$ws.Cells.Item(50, 5).Value = '  -1.40%  '

$ws.Cells.Item(51, 4).Value = '0.3404'
$ws.Cells.Item(51, 5).Value = '  -0.86%  '
